# Add the "Marisa Trie #4 / CDS space" column (P) results for the cptVSsd_cpt
# comparison sheet, matching the analogous existing Marisa Trie columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Fill in the new P2:P11 data values (copy the cell format from a
#        neighbouring "Marisa Trie / CDS space" column that already carries
#        the right style, then overwrite the value). ---

$ws.Range("N2").Copy()
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("P2").Value = 0.117176

$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("P3").Value = 0.15879199999999999

$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 0.158192

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value = 0.0066160000000000004

$ws.Range("N6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = 0.08856

$ws.Range("O7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P7").Value = 0.01992

$ws.Range("O8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value = 0.088496000000000005

$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P9").Value = 0.049784000000000002

$ws.Range("O10").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("P10").Value = 0.17586399999999999

$ws.Range("O11").Copy()
$ws.Range("P11").PasteSpecial(-4122)
$ws.Range("P11").Value = 0.209616

$excel.CutCopyMode = $false

# --- 2. Re-enter the ratio formulas for P19:P28 (same pattern as the other
#        ratio columns L:O, one row per dataset) so they recompute against
#        the newly-populated P2:P11. Each cell is written individually
#        (rather than via AutoFill/PasteSpecial) so it correctly joins the
#        existing shared-formula group anchored on the neighbouring O column
#        for that row. ---

$ws.Range("P19").Formula = "=P2/`$F2"
$ws.Range("P20").Formula = "=P3/`$F3"
$ws.Range("P21").Formula = "=P4/`$F4"
$ws.Range("P22").Formula = "=P5/`$F5"
$ws.Range("P23").Formula = "=P6/`$F6"
$ws.Range("P24").Formula = "=P7/`$F7"
$ws.Range("P25").Formula = "=P8/`$F8"
$ws.Range("P26").Formula = "=P9/`$F9"
$ws.Range("P27").Formula = "=P10/`$F10"
$ws.Range("P28").Formula = "=P11/`$F11"

# --- 3. The scratch/staging cells that used to hold placeholder P values
#        further down the sheet (P33:P42) are no longer needed now that the
#        real data lives in P2:P11 - remove them entirely (format + content)
#        so the sheet's used range shrinks back down. ---

$ws.Range("P33:P42").Clear()

# --- 4. Restore the view state (scrolled right to show the new column,
#        selection resting on the newly filled formula range). ---

$ws.Range("P19:P28").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1

$wb.Application.CalculateFull()
